$d = $word.ActiveDocument

# 1. "How much time outside the class, on average, do you think you spent?"
#    -> "On average, how much time outside the class do you think you spent on CS-related work?"
$d.Content.Find.Execute(
    "How much time outside the class, on average, do you think you spent?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "On average, how much time outside the class do you think you spent on CS-related work?",
    2)

# 2. "mountains and tree" (followed by a separate run "s.")
#    -> "mountains and fractal tree" (so the sentence reads "... fractal mountains and fractal trees.")
$d.Content.Find.Execute(
    "mountains and tree",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mountains and fractal tree",
    2)
